$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill previously-empty cells that now carry values (processing of empty cells)
$ws.Range("C3").Value = 45646
$ws.Range("D3").Value = "m"
$ws.Range("D4").Value = "k"
$ws.Range("D5").Value = "k"
$ws.Range("C10").Value = 5
$ws.Range("C11").Value = 546

# Update the view state: scrolled so row 4 is at top, selection on D10
$ws.Range("D10").Select()
$excel.ActiveWindow.ScrollRow = 4
